$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "69.385.15"
$ws.Range("E2").Value = "  +2.88%  "

$ws.Range("D3").Value = "3.763.75"
$ws.Range("E3").Value = "  +1.66%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.11%  "

Set-TextValue $ws.Range("D5") "605.14"
$ws.Range("E5").Value = "  +1.55%  "

Set-TextValue $ws.Range("D6") "169.73"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("D7").Value = "3.759.82"
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("E8").Value = "  -0.05%  "

Set-TextValue $ws.Range("D9") "0.539"
$ws.Range("E9").Value = "  +2.98%  "

Set-TextValue $ws.Range("D10") "0.169"
$ws.Range("E10").Value = "  +5.34%  "

Set-TextValue $ws.Range("D11") "6.40"
$ws.Range("E11").Value = "  +3.22%  "

Set-TextValue $ws.Range("D12") "0.464"
$ws.Range("E12").Value = "  +0.45%  "

Set-TextValue $ws.Range("D13") "38.47"
$ws.Range("E13").Value = "  +2.40%  "

$ws.Range("D15").Value = "4.382.24"
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("D16").Value = "3.756.84"
$ws.Range("E16").Value = "  +1.63%  "

$ws.Range("D17").Value = "69.242.32"
$ws.Range("E17").Value = "  +2.55%  "

Set-TextValue $ws.Range("D18") "7.32"
$ws.Range("E18").Value = "  +2.10%  "

$ws.Range("E19").Value = "  -0.10%  "

Set-TextValue $ws.Range("D20") "17.14"
$ws.Range("E20").Value = "  -2.57%  "

Set-TextValue $ws.Range("D21") "10.90"
$ws.Range("E21").Value = "  +19.06%  "

Set-TextValue $ws.Range("D22") "496.38"
$ws.Range("E22").Value = "  +0.96%  "

Set-TextValue $ws.Range("D23") "0.730"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("E24").Value = "  +11.92%  "

Set-TextValue $ws.Range("D25") "85.50"
$ws.Range("E25").Value = "  -0.47%  "

Set-TextValue $ws.Range("D26") "2.34"
$ws.Range("E26").Value = "  +1.67%  "

Set-TextValue $ws.Range("D27") "12.40"
$ws.Range("E27").Value = "  +1.82%  "

Set-TextValue $ws.Range("D28") "10.31"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("E29").Value = "  -0.20%  "

Set-TextValue $ws.Range("D30") "2.53"
$ws.Range("E30").Value = "  +7.53%  "

$ws.Range("E31").Value = "  +2.07%  "

Set-TextValue $ws.Range("D32") "7.96"
$ws.Range("E32").Value = "  +4.29%  "

Set-TextValue $ws.Range("D33") "32.09"
$ws.Range("E33").Value = "  +1.44%  "

$ws.Range("D34").Value = "3.900.77"
$ws.Range("E34").Value = "  +1.61%  "

$ws.Range("E35").Value = "  +1.48%  "

$ws.Range("D36").Value = "3.689.02"
$ws.Range("E36").Value = "  +1.29%  "

Set-TextValue $ws.Range("D37") "0.998"
$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("E38").Value = "  +2.00%  "

Set-TextValue $ws.Range("D39") "5.91"
$ws.Range("E39").Value = "  +2.78%  "

$ws.Range("E40").Value = "  +1.72%  "

Set-TextValue $ws.Range("D41") "0.325"
$ws.Range("E41").Value = "  +0.80%  "

Set-TextValue $ws.Range("D42") "3.02"
$ws.Range("E42").Value = "  +8.67%  "

Set-TextValue $ws.Range("D43") "438.02"
$ws.Range("E43").Value = "  +0.75%  "

Set-TextValue $ws.Range("D44") "48.63"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("E45").Value = "  +2.88%  "

Set-TextValue $ws.Range("D46") "8.49"
$ws.Range("E46").Value = "  +1.36%  "

Set-TextValue $ws.Range("D48") "40.50"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").Value = "2.813.54"
$ws.Range("E49").Value = "  +2.08%  "

Set-TextValue $ws.Range("D50") "141.16"
$ws.Range("E50").Value = "  -1.12%  "

$ws.Range("E51").Value = "  +2.58%  "
